$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '63.310.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +5.61%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '3.382.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  +5.99%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  -0.04%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '574.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +7.20%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'" + '154.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + '  +6.02%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  -0.11%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'" + '3.382.26'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + '  +5.78%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + '0.529'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  -0.33%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'" + '7.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  +2.01%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'" + '0.120'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  +6.99%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + '0.437'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  +1.72%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '3.960.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  +5.86%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +0.18%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'" + '0.0000184'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  +6.58%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '27.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  +5.04%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '63.326.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  +5.60%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'" + '3.364.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  +4.71%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + '6.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  +1.67%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '13.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  +5.06%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + '8.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  +2.87%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '388.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  +5.11%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  -0.04%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'" + '0.538'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'" + '  +3.13%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + '70.70'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  +1.88%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '9.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  +12.08%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +6.47%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'" + '0.0000102'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'" + '  +17.53%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'" + '  -0.38%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '2.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  +7.39%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'" + 'RenderToken'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'" + '6.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  +6.01%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'" + 'EthereumClassic'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'" + 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'" + '23.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + '  +3.41%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'" + 'NEARProtocol'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'" + 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'" + '5.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  +6.28%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '1.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  +9.75%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + '6.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +3.43%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '1.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +10.25%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + '158.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  +1.05%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'" + 'EnergySwap'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'" + '27.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  +5.68%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'" + 'Stacks'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'" + 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'" + '1.88'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +12.67%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '2.911.03'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  +3.05%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '0.0748'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +6.21%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + '0.0329'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  +8.46%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'" + '41.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + '  +2.95%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'" + '0.754'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  +5.28%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '4.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  +1.54%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + '1.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  +6.61%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '3.423.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  +5.92%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'" + 'InjectiveProtocol'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'" + 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'" + '22.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  +6.80%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'" + 'Bittensor'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'" + 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'" + '300.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  +13.63%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'" + '0.104'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  -1.07%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + '6.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  +2.62%  '
$ws.Range("E51").Style = "Normal"
